$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data rows (column B only) below the existing table.
$ws.Range("B17").Value = 50
$ws.Range("B18").Value = 500
$ws.Range("B19").Value = 1000
$ws.Range("B20").Value = 3000
$ws.Range("B21").Value = 5000
$ws.Range("B22").Value = 7000
$ws.Range("B23").Value = 10000

# Reposition/resize the chart (dragged up and to the left, made a bit
# narrower and slightly taller). Values are in points, matching the
# EMU anchor recorded in the drawing XML (12700 EMU per point).
$co = $ws.ChartObjects(1)
$co.Left = 330.4306640625
$co.Top = 3
$co.Width = 192.3125
$co.Height = 178

# Move the active selection to match where the user left the cursor.
$ws.Range("C17").Select()
